$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text run formatting) ---
$ws.Range("A8").Characters(21, 2).Text = "23"
$ws.Range("A8").Characters(1, 7).Font.Size = 10
$ws.Range("A8").Characters(8, 2).Font.Size = 10
$ws.Range("A8").Characters(10, 11).Font.Size = 10
$ws.Range("A8").Characters(21, 2).Font.Size = 10

$ws.Range("C9").Characters(27, 9).Text = "6/3/2024"
$ws.Range("C9").Characters(46, 8).Text = "6/9/2024"
$ws.Range("C9").Characters(1, 26).Font.Size = 10
$ws.Range("C9").Characters(27, 8).Font.Size = 10
$ws.Range("C9").Characters(35, 11).Font.Size = 10
$ws.Range("C9").Characters(46, 8).Font.Size = 10

# --- Crime statistics table updates ---
# Row 15
$ws.Range("G15").Value = 1
$ws.Range("L15").Value = 50

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = -3.703703703703
$ws.Range("L16").Value = -31.578947368421
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -82.666666666666

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 600
$ws.Range("F17").Value = 15
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 44
$ws.Range("J17").Value = 46
$ws.Range("K17").Value = -4.347826086956
$ws.Range("L17").Value = 37.5
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -18.518518518518

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 250
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 13.043478260869
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 150
$ws.Range("K18").Value = -22.666666666666
$ws.Range("L18").Value = -4.132231404958
$ws.Range("M18").Value = 1.754385964912
$ws.Range("N18").Value = -75.160599571734

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 120
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 10.714285714285
$ws.Range("I19").Value = 207
$ws.Range("J19").Value = 296
$ws.Range("K19").Value = -30.067567567567
$ws.Range("L19").Value = -15.853658536585
$ws.Range("M19").Value = 29.375
$ws.Range("N19").Value = -8.407079646017

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 54.545454545454
$ws.Range("I20").Value = 114
$ws.Range("J20").Value = 65
$ws.Range("K20").Value = 75.384615384615
$ws.Range("L20").Value = 178.048780487805
$ws.Range("M20").Value = 83.870967741935
$ws.Range("N20").Value = -91.920623671155

# Row 21
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 130.769230769231
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = 18.987341772151
$ws.Range("I21").Value = 510
$ws.Range("J21").Value = 590
$ws.Range("K21").Value = -13.559322033898
$ws.Range("L21").Value = 6.25
$ws.Range("M21").Value = 27.5
$ws.Range("N21").Value = -77.969762419006

# Row 24
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -44.444444444444
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = 15.78947368421
$ws.Range("I24").Value = 245
$ws.Range("J24").Value = 235
$ws.Range("K24").Value = 4.255319148936
$ws.Range("L24").Value = -30.397727272727
$ws.Range("M24").Value = 17.224880382775

# Row 25
$ws.Range("C25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 49
$ws.Range("J25").Value = 49
$ws.Range("L25").Value = -2

# Row 26
$ws.Range("C26").Value = 3
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 86
$ws.Range("J26").Value = 99
$ws.Range("K26").Value = -13.131313131313
$ws.Range("L26").Value = -17.307692307692
$ws.Range("M26").Value = 22.857142857142

# Row 27
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = 66.666666666666

# Row 28
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = -66.666666666666

# --- Cells converted from numeric to shared-text "0"/" " placeholders ---
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("C28").Copy($ws.Range("D28"))
$ws.Range("M28").Copy($ws.Range("E28"))
